$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append the new "Spillkråka" section at the end of the body,
#    right after the last paragraph ("SLU Artdatabanken, Uppsala ")
#    and before the sectPr.
# ------------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)

# -- Heading1: "Spillkråka – ekologi samt krav på livsmiljön"
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p1.Range.Text = "Spillkråka – ekologi samt krav på livsmiljön"
$p1.Range.Font.Italic = $False
$p1.Style = "Heading1"

# -- Normal body paragraph about spillkråka ecology
$p1.Range.Collapse(0)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2.Range.Text = "Spillkråka (NT) är rödlistad som nära hotad, fridlyst enligt §4 Artskyddsförordningen och ingår i bilaga 1 i EU:s fågeldirektiv. Spillkråka lever i både barr- och blandskog liksom i ren lövskog. De tätaste populationerna tenderar att finnas i äldre, variationsrik blandskog med gott om död ved och gamla träd. Varje par utnyttjar 400-1 000 hektar skog beroende på skogens kvalitet. En minskning av populationen pågår på grund av minskad tillgång på lämpliga bo- och födoträd och minskad födotillgång. Spillkråkans minskningstakt har uppgått till 19 (24-10) % under de senaste 15 åren. Skogsbruk med korta omloppstider och täta, homogena ungskogar utgör det största hotet (Artdatabanken 2023)."
$p2.Range.Font.Italic = $False
$p2.Style = "Normal"

# -- Heading2: "Referenser - spillkråka"
$p2.Range.Collapse(0)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p3.Range.Text = "Referenser - spillkråka"
$p3.Range.Font.Italic = $False
$p3.Style = "Heading2"

# -- Reference paragraph with 3 runs: plain / italic / plain
$p3.Range.Collapse(0)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$run1 = "SLU Artdatabanken, 2021. "
$run2 = "Artfaktablad. Naturvård – artfakta. "
$run3 = "SLU Artdatabanken, Uppsala"
$p4.Range.Text = $run1 + $run2 + $run3
$p4.Range.Font.Italic = $False
$p4.Style = "Normal"

$p4start = $p4.Range.Start
$italicStart = $p4start + $run1.Length
$italicEnd = $italicStart + $run2.Length
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Font.Italic = $True

# ------------------------------------------------------------------
# 2) Update the date shown in the first-page header
#    (wdHeaderFooterFirstPage = 2) from 2023-10-13 to 2023-10-22.
# ------------------------------------------------------------------

$header = $d.Sections(1).Headers(2)
$header.Range.Find.Execute("2023-10-13", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "2023-10-22", 2)
